$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# New "CA Images" picture links added alongside the Horizontal Cabling and
# Pathways section (rows 30-36), referencing shared strings 53-58.
$ws.Range("F30").Value = "CA Images/Picture1"
$ws.Range("H30").Value = "CA Images/Picture 2"
$ws.Range("F31").Value = "CA Images/Picture 3"
$ws.Range("F33").Value = "CA Images/Picture 4"
$ws.Range("F34").Value = "CA Images/Picture 5"
$ws.Range("F36").Value = "CA Images/Picture 6"

# Update the active selection to match where the user ended up working.
$ws.Range("F30").Select()
